# The commit inserts a new bold/underlined "MULTI PARADIGM " run right at
# the start of the heading paragraph (before the existing bookmark /
# "OF PROGRAMMING LANGUAGE..." run), matching that run's character
# formatting (rFonts cstheme=minorHAnsi, bold, sz/szCs 24, single underline).

$d = $word.ActiveDocument

# Second paragraph holds the heading text
# "OF PROGRAMMING LANGUAGE; Python Programming Language."
$heading = $d.Paragraphs(2)

# Collapsed range at the very start of that paragraph (right before the
# bookmark / first run) so the new run is inserted as its own element.
$insertAt = $heading.Range.Start
$r = $d.Range($insertAt, $insertAt)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:u w:val="single"/>
              </w:rPr>
              <w:t xml:space="preserve">MULTI PARADIGM </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
